# Apply the "Add files via upload" revision to analisis1.xlsx
#
# Summary of the edit (reconstructed from the target OOXML diff):
#  1. Columns E:F (an empty spacer column plus a now-unused helper column)
#     are deleted, shifting the L:O helper columns left to J:M.
#  2. Two of the ticker labels that are no longer used ("ETTB", "ETB E")
#     are consolidated: row 24's ticker becomes "ETB" (matching the
#     existing padded "ETB       " label used elsewhere) and row 30's
#     ticker becomes "ETB " (a new, shorter label).
#  3. The leftover "PARCIAL" grading scratch area (rows 49-52) has its
#     contents cleared: row 49 (the "PARCIAL" caption) disappears
#     entirely, and the formulas/values in C50:D50 / C52:D52 are cleared
#     while keeping their number formats; the extra E/F cells in that
#     block go away together with the deleted columns.
#  4. The active selection moves to B31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the two now-unused columns (E:F), shifting everything right
#    of them two columns to the left.
$ws.Range("E1:F1").EntireColumn.Delete()

# 2. Relabel the two duplicate/obsolete "ETB"-family tickers.
$ws.Range("B24").Value = "ETB       "
$ws.Range("B30").Value = "ETB "

# 3. Clear out the old "PARCIAL" grade scratch pad, keeping the styled
#    but now-empty cells in rows 50 and 52.
$ws.Range("C49").ClearContents()
$ws.Range("C50:D50").ClearContents()
$ws.Range("C52:D52").ClearContents()

# 4. Move the selection to B31, matching the saved workbook state.
$ws.Range("B31").Select()
